$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.03"
$ws.Range("E2").Value = "'3.19%"
$ws.Range("D3").Value = "'35.40"
$ws.Range("E3").Value = "'-0.26%"
$ws.Range("D4").Value = "'5.101"
$ws.Range("E4").Value = "'0.83%"
$ws.Range("D5").Value = "'0.08166"
$ws.Range("E5").Value = "'3.66%"
$ws.Range("D6").Value = "'2.085"
$ws.Range("E6").Value = "'-2.49%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.939"
$ws.Range("E7").Value = "'-0.07%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9314"
$ws.Range("E8").Value = "'0.98%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1032"
$ws.Range("E9").Value = "'6.49%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1929"
$ws.Range("E10").Value = "'4.77%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09103"
$ws.Range("E11").Value = "'4.53%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03589"
$ws.Range("E12").Value = "'0.53%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09912"
$ws.Range("E13").Value = "'0.15%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001443"
$ws.Range("E14").Value = "'0.24%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005690"
$ws.Range("E15").Value = "'0.92%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.469"
$ws.Range("E16").Value = "'-0.03%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.137"
$ws.Range("E17").Value = "'0.36%"
$ws.Range("D18").Value = "'2.853"
$ws.Range("E18").Value = "'3.63%"
$ws.Range("E19").Value = "'2.53%"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'-1.36%"
$ws.Range("D21").Value = "'5.098"
$ws.Range("E21").Value = "'-1.20%"
$ws.Range("D22").Value = "'0.2215"
$ws.Range("E22").Value = "'0.34%"
$ws.Range("D23").Value = "'0.04552"
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("D24").Value = "'0.001242"
$ws.Range("E24").Value = "'0.85%"
$ws.Range("D25").Value = "'0.004795"
$ws.Range("E25").Value = "'-0.73%"
$ws.Range("E26").Value = "'-3.74%"
$ws.Range("D27").Value = "'0.0004506"
$ws.Range("E27").Value = "'-5.11%"
$ws.Range("D39").Value = "'0.01978"
$ws.Range("E39").Value = "'6.92%"
$ws.Range("D40").Value = "'0.04919"
$ws.Range("E40").Value = "'3.97%"
$ws.Range("D41").Value = "'0.007610"
$ws.Range("E41").Value = "'-2.19%"
$ws.Range("D42").Value = "'0.1385"
$ws.Range("E42").Value = "'0.04%"
$ws.Range("D43").Value = "'0.007877"
$ws.Range("E43").Value = "'1.80%"
$ws.Range("D44").Value = "'0.002188"
$ws.Range("E44").Value = "'1.16%"
$ws.Range("D45").Value = "'0.01173"
$ws.Range("E45").Value = "'3.59%"
$ws.Range("D46").Value = "'0.00006683"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("D48").Value = "'192.10"
$ws.Range("E48").Value = "'279.79%"
$ws.Range("D49").Value = "'0.001702"
$ws.Range("E49").Value = "'-10.40%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.10%"
